$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the Password values between row 2 and row 3 (B2 <-> B3)
$b2 = $ws.Range("B2").Value2
$b3 = $ws.Range("B3").Value2
$ws.Range("B2").Value = $b3
$ws.Range("B3").Value = $b2

# Move the active selection to A7
$ws.Range("A7").Select() | Out-Null
